$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("G2").Value = 17.040132
$ws.Range("H2").Value = 51.12039599999999
$ws.Range("I2").Value = 0.3748878535193673
$ws.Range("J2").Value = 0.3748878535193674
$ws.Range("M2").Value = 113.5893336666667
$ws.Range("N2").Value = 340.768001
$ws.Range("O2").Value = 0.8306211829777892
$ws.Range("P2").Value = 0.8306211829777892
$ws.Range("Q2").Value = 1935.577239472044
$ws.Range("R2").Value = 17420.1951552484
$ws.Range("S2").Value = 0.311389792374261
$ws.Range("T2").Value = 0.3113897923742611

# Row 3
$ws.Range("G3").Value = 17.040132
$ws.Range("H3").Value = 51.12039599999999
$ws.Range("I3").Value = 0.3748878535193673
$ws.Range("J3").Value = 0.3748878535193674
$ws.Range("O3").Value = 0.1531387528565491
$ws.Range("P3").Value = 0.1531387528565491
$ws.Range("Q3").Value = 356.8556769135479
$ws.Range("R3").Value = 3211.701092221931
$ws.Range("S3").Value = 0.05740985834902458
$ws.Range("T3").Value = 0.05740985834902458

# Row 4
$ws.Range("G4").Value = 17.040132
$ws.Range("H4").Value = 51.12039599999999
$ws.Range("I4").Value = 0.3748878535193673
$ws.Range("J4").Value = 0.3748878535193674
$ws.Range("M4").Value = 2.220865666666667
$ws.Range("N4").Value = 6.662597
$ws.Range("O4").Value = 0.01624006416566169
$ws.Range("P4").Value = 0.01624006416566169
$ws.Range("Q4").Value = 37.84384411426799
$ws.Range("R4").Value = 340.594597028412
$ws.Range("S4").Value = 0.006088202796081708
$ws.Range("T4").Value = 0.006088202796081708

# Row 5
$ws.Range("I5").Value = 0.1732576994636535
$ws.Range("J5").Value = 0.1732576994636535
$ws.Range("M5").Value = 113.5893336666667
$ws.Range("N5").Value = 340.768001
$ws.Range("O5").Value = 0.8306211829777892
$ws.Range("P5").Value = 0.8306211829777892
$ws.Range("Q5").Value = 894.543945601082
$ws.Range("R5").Value = 8050.895510409738
$ws.Range("S5").Value = 0.1439115152885101
$ws.Range("T5").Value = 0.1439115152885101

# Row 6
$ws.Range("I6").Value = 0.1732576994636535
$ws.Range("J6").Value = 0.1732576994636535
$ws.Range("O6").Value = 0.1531387528565491
$ws.Range("P6").Value = 0.1531387528565491
$ws.Range("S6").Value = 0.02653246801865869
$ws.Range("T6").Value = 0.02653246801865869

# Row 7
$ws.Range("I7").Value = 0.1732576994636535
$ws.Range("J7").Value = 0.1732576994636535
$ws.Range("M7").Value = 2.220865666666667
$ws.Range("N7").Value = 6.662597
$ws.Range("O7").Value = 0.01624006416566169
$ws.Range("P7").Value = 0.01624006416566169
$ws.Range("Q7").Value = 17.489863457954
$ws.Range("R7").Value = 157.408771121586
$ws.Range("S7").Value = 0.002813716156484662
$ws.Range("T7").Value = 0.002813716156484662

# Row 8
$ws.Range("G8").Value = 20.53856733333333
$ws.Range("H8").Value = 61.615702
$ws.Range("I8").Value = 0.4518544470169792
$ws.Range("J8").Value = 0.4518544470169792
$ws.Range("M8").Value = 113.5893336666667
$ws.Range("N8").Value = 340.768001
$ws.Range("O8").Value = 0.8306211829777892
$ws.Range("P8").Value = 0.8306211829777892
$ws.Range("Q8").Value = 2332.9621778613
$ws.Range("R8").Value = 20996.6596007517
$ws.Range("S8").Value = 0.375319875315018
$ws.Range("T8").Value = 0.375319875315018

# Row 9
$ws.Range("G9").Value = 20.53856733333333
$ws.Range("H9").Value = 61.615702
$ws.Range("I9").Value = 0.4518544470169792
$ws.Range("J9").Value = 0.4518544470169792
$ws.Range("O9").Value = 0.1531387528565491
$ws.Range("P9").Value = 0.1531387528565491
$ws.Range("Q9").Value = 430.1201627177037
$ws.Range("R9").Value = 3871.081464459333
$ws.Range("S9").Value = 0.06919642648886584
$ws.Range("T9").Value = 0.06919642648886583

# Row 10
$ws.Range("G10").Value = 20.53856733333333
$ws.Range("H10").Value = 61.615702
$ws.Range("I10").Value = 0.4518544470169792
$ws.Range("J10").Value = 0.4518544470169792
$ws.Range("M10").Value = 2.220865666666667
$ws.Range("N10").Value = 6.662597
$ws.Range("O10").Value = 0.01624006416566169
$ws.Range("P10").Value = 0.01624006416566169
$ws.Range("Q10").Value = 45.61339903312155
$ws.Range("R10").Value = 410.520591298094
$ws.Range("S10").Value = 0.007338145213095324
$ws.Range("T10").Value = 0.007338145213095324
